# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Updates Price (D) and Volume(1h) (E) values for most rows, and swaps the
# Chainlink/BitcoinCash (rows 20-21) and USDe/Monero (rows 48-49) rankings.
#
# Price values that look numeric (e.g. "612.91") must be written as text to
# match the source data (prices like "69.508.02" or "2.950.34" use '.' as a
# thousands separator and are never real Excel numbers), so for those cells
# we force the NumberFormat to "@" (Text) before assigning the string, then
# restore the cell style to "Normal" so no stray number formatting lingers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.508.02"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "3.729.08"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("D7").Value = "3.727.58"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.57"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.94%  "
$ws.Range("E12").Value = "  -3.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").Value = "4.349.85"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "3.727.40"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "69.626.20"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "501.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.71%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.99%  "
$ws.Range("E23").Value = "  -1.10%  "
$ws.Range("E24").Value = "  +5.15%  "
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.02%  "
$ws.Range("E28").Value = "  +7.56%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.351"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.13%  "
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("E40").Value = "  +12.91%  "
$ws.Range("E41").Value = "  -5.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "439.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "49.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("E45").Value = "  -3.27%  "
$ws.Range("D46").Value = "2.950.34"
$ws.Range("E46").Value = "  -4.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0359"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "138.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("E51").Value = "  -1.58%  "
